# Contract_WIP.xlsx - append 3 new PO rows to the bottom of the WIP list
# (rows 30-32), matching the style of the existing "text" rows further up
# the sheet (e.g. rows 9+): every field is stored as plain text, including
# the date-looking and currency-looking values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 30

$newRows = @(
    @("4/1/2019", "SPE7M5-19-V-7228", "1",    "`$629.00",    "5935015527394", "BACKSHELL,ELECTRICAL CONNECTOR", "Glenair", "390AS002NF0804H3", "CP", "2019 SEP 09"),
    @("4/1/2019", "SPE4A6-19-V-070V", "14",   "`$17,032.68", "5342015037327", "MAGNET,HOLDER ASSY",             "GEMS",    "137444",            "ZZ", "2019 SEP 09"),
    @("4/2/2019", "SPE7M5-19-P-6235", "1500", "`$65,145.00", "5935016786944", "BACKSHELL,ELECTRICAL CONNECTOR", "Glenair", "445HS065NF25064",   "41", "2019 SEP 19")
)

$endRow = $startRow + $newRows.Count - 1

# Pre-format the target block as Text so values like "4/1/2019" or
# "$629.00" land as literal strings instead of being auto-converted to
# dates / currency numbers by Excel's input parsing.
$targetRange = $ws.Range("A" + $startRow + ":J" + $endRow)
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        $cell.Value = $rowData[$col - 1]
    }
}
